# 03_Edit_or_Delete_ComputerEntry_Page_Validations.xlsx
# Adds a new validation row-note to the "Edit and Delete" CRUD test-case sheet:
# every step row (rows 8-22, column B) gets the note that the home page should
# have records and the user clicked a computer name link in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$note = "Home page should have records and the user clicked on a computer name link in the table"

# Fill column B for rows 8 through 22 with the new note text.
for ($r = 8; $r -le 22; $r++) {
    $ws.Cells.Item($r, 2).Value = $note
}

# A handful of rows need to grow taller to fit the extra wrapped text.
$tallerRows = @(9, 10, 11, 12, 14, 16, 18)
foreach ($r in $tallerRows) {
    $ws.Rows.Item($r).RowHeight = 45
}

# Leave the view scrolled down to, and focused on, the last edited cell.
$null = $ws.Range("C22").Select()
